# Oracle Apps Course Content added in IT Syllabus file
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the Tally sheet: bold the "Link" label cell (A20)
# ---------------------------------------------------------------------------
$wsTally = $wb.Worksheets.Item("Tally")
$wsTally.Range("A20").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Add the new "OracleApps" worksheet as the last tab in the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "OracleApps"

# Numbered topic list
$topics = @(
    "Introduction to Oracle Apps",
    "Oracle Applications",
    "Product Directory Structure",
    "Data Model",
    "Responsibilities & User",
    "New Module Development",
    "Table Registration",
    "FlexFields",
    "New Form Development",
    "Defining Calendars",
    "WHO Columns",
    "Non-Form Functions",
    "Search Methods",
    "Profiles",
    "Customization of Forms",
    "Concurrent Processing(CP)",
    "Flex Field Reports",
    "Qualifiers",
    "Interfaces",
    "Business Components",
    "Multi Organizations",
    "Alerts",
    "Discoverer"
)

$row = 4
for ($i = 0; $i -lt $topics.Length; $i++) {
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $topics[$i]
    $row++
}

# Header
$ws.Range("B2").Value = "Oracle Apps Topics"
$ws.Range("B2").Font.Bold = $true

# Note row
$ws.Range("A29").Value = "Note"
$ws.Range("B29").Value = "Needs to verify above contents of course from Oracle Apps expert professional"
$ws.Range("A29:D29").Font.Bold = $true

$ws.Columns("B:B").ColumnWidth = 29.0

$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. Restore the selections / active-tab state so it matches the saved file
# ---------------------------------------------------------------------------
$wsTally.Activate()
$wsTally.Range("A20").Select() | Out-Null

$ws.Activate()
$ws.Range("E3").Select() | Out-Null
